# Applies the "Updated symbol list" data refresh for rows 2-51:
#   - Column D (Price): new quote for each coin
#   - Column E (Volume(1h)): new 1h change percentage
#   - Column G (Hora): updated run hour, 20 -> 21, for every row
# All three columns hold plain text in this sheet (no numeric
# values), so the Text number format ("@") is applied before each
# write to stop Excel from reinterpreting the strings as numbers,
# percentages, dates, etc.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column D (Price) -------------------------------------------------
$priceUpdates = @{
    2 = "310.07"
    3 = "40.69"
    4 = "5.097"
    5 = "0.07837"
    6 = "4.345"
    7 = "1.675"
    8 = "0.9204"
    9 = "0.1061"
    10 = "0.1776"
    11 = "0.09011"
    12 = "0.04404"
    13 = "7.216"
    15 = "0.001289"
    16 = "0.006028"
    17 = "3.380"
    18 = "2.575"
    19 = "0.3368"
    20 = "0.1386"
    21 = "0.2805"
    22 = "0.04178"
    23 = "0.001221"
    24 = "0.004139"
    25 = "0.0001228"
    26 = "0.0002998"
    38 = "0.02434"
    39 = "0.05262"
    40 = "0.008002"
    41 = "0.1353"
    42 = "0.007588"
    43 = "0.001994"
    44 = "0.008072"
    45 = "0.3378"
    46 = "0.00006778"
    47 = "0.00000000755"
    48 = "0.003421"
    49 = "0.004127"
    50 = "0.00002114"
    51 = "0.0002013"
}
foreach ($row in $priceUpdates.Keys) {
    $cell = $ws.Cells.Item($row, 4)   # column D
    $cell.NumberFormat = "@"
    $cell.Value = $priceUpdates[$row]
}

# --- Column E (Volume(1h)) ---------------------------------------------
$volumeUpdates = @{
    2 = "-4.99%"
    3 = "-7.86%"
    4 = "-3.40%"
    5 = "-5.88%"
    6 = "-1.23%"
    7 = "-13.49%"
    8 = "-4.83%"
    9 = "-5.68%"
    10 = "-6.00%"
    11 = "-7.25%"
    12 = "-4.14%"
    13 = "-15.80%"
    14 = "-0.05%"
    15 = "0.05%"
    16 = "2.75%"
    17 = "-0.32%"
    18 = "1.02%"
    19 = "0.30%"
    20 = "1.23%"
    21 = "8.68%"
    22 = "0.57%"
    23 = "-1.06%"
    24 = "-6.06%"
    25 = "-5.78%"
    26 = "0.59%"
    38 = "-9.33%"
    39 = "-5.19%"
    40 = "2.30%"
    41 = "-3.73%"
    42 = "3.67%"
    43 = "-5.68%"
    44 = "2.55%"
    45 = "-3.61%"
    46 = "-0.91%"
    47 = "0.58%"
    48 = "-1.88%"
    49 = "16.83%"
    50 = "0.58%"
    51 = "0.58%"
}
foreach ($row in $volumeUpdates.Keys) {
    $cell = $ws.Cells.Item($row, 5)   # column E
    $cell.NumberFormat = "@"
    $cell.Value = $volumeUpdates[$row]
}

# --- Column G (Hora) ----------------------------------------------------
# Every data row (2-51) moves from hour "20" to hour "21".
$horaRange = $ws.Range("G2:G51")
$horaRange.NumberFormat = "@"
$horaRange.Value = "21"

